$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 346.75
$ws.Range("I2").Value = 337.25
$ws.Range("J2").Value = 356.25
$ws.Range("K2").Value = 337.25
$ws.Range("L2").Value = 356.25
$ws.Range("M2").Value = -224.25
$ws.Range("N2").Value = -582.25
$ws.Range("H41").Value = 1722.1111
$ws.Range("I41").Value = 2650.25
$ws.Range("J41").Value = 979.6
$ws.Range("K41").Value = 2650.25
$ws.Range("L41").Value = 979.6
$ws.Range("M41").Value = -2210.25
$ws.Range("N41").Value = -1859.6
$ws.Range("H53").Value = 483.66666
$ws.Range("I53").Value = 484.6
$ws.Range("J53").Value = 479
$ws.Range("K53").Value = 484.6
$ws.Range("L53").Value = 479
$ws.Range("M53").Value = 152.4
$ws.Range("N53").Value = -1753
$ws.Range("H58").Value = 665.8
$ws.Range("I58").Value = 665.8
$ws.Range("K58").Value = 1997.4
$ws.Range("M58").Value = -1847.4
$ws.Range("H70").Value = 13038.5
$ws.Range("J70").Value = 16884.666
$ws.Range("L70").Value = 50653.99800000001
$ws.Range("N70").Value = -51193.99800000001
$ws.Range("H73").Value = 13038.5
$ws.Range("J73").Value = 16884.666
$ws.Range("L73").Value = 50653.99800000001
$ws.Range("N73").Value = -52525.99800000001
$ws.Range("H98").Value = 2392.818
$ws.Range("I98").Value = 2021.625
$ws.Range("K98").Value = 2021.625
$ws.Range("M98").Value = -523.625
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H113").Value = 4000
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 254
$ws.Range("N113").Value = -11508
$ws.Range("H116").Value = 2950
$ws.Range("I116").Value = 2900
$ws.Range("K116").Value = 2900
$ws.Range("M116").Value = 542
$ws.Range("H122").Value = 2392.818
$ws.Range("I122").Value = 2021.625
$ws.Range("K122").Value = 6064.875
$ws.Range("M122").Value = -3614.875
$ws.Range("H132").Value = 3756.4285
$ws.Range("I132").Value = 3116.0588
$ws.Range("J132").Value = 6478
$ws.Range("K132").Value = 9348.1764
$ws.Range("L132").Value = 19434
$ws.Range("M132").Value = -6818.1764
$ws.Range("N132").Value = -24494

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 510

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 510
$ws.Range("H12").Value = 6832.6665
$ws.Range("I12").Value = 9999
$ws.Range("K12").Value = 9999
$ws.Range("M12").Value = -9831
$ws.Range("H105").Value = 3403
$ws.Range("I105").Value = 3477.5
$ws.Range("J105").Value = 3303.6667
$ws.Range("K105").Value = 3477.5
$ws.Range("L105").Value = 3303.6667
$ws.Range("M105").Value = -1730.5
$ws.Range("N105").Value = -6797.6667
$ws.Range("H107").Value = 1915.6666
$ws.Range("I107").Value = 1948
$ws.Range("J107").Value = 1883.3334
$ws.Range("K107").Value = 1948
$ws.Range("L107").Value = 1883.3334
$ws.Range("M107").Value = -28
$ws.Range("N107").Value = -5723.3334

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 90.46154
$ws.Range("I7").Value = 68.125
$ws.Range("J7").Value = 126.2
$ws.Range("K7").Value = 68.125
$ws.Range("L7").Value = 126.2
$ws.Range("M7").Value = 44.875
$ws.Range("N7").Value = -352.2
$ws.Range("H17").Value = 1200
$ws.Range("I17").Value = 1200
$ws.Range("K17").Value = 1200
$ws.Range("M17").Value = -1026
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H56").Value = 6250
$ws.Range("I56").Value = 6250
$ws.Range("K56").Value = 6250
$ws.Range("M56").Value = -5405
$ws.Range("H105").Value = 1883.1666
$ws.Range("I105").Value = 1612.25
$ws.Range("J105").Value = 2425
$ws.Range("K105").Value = 1612.25
$ws.Range("L105").Value = 2425
$ws.Range("M105").Value = 134.75
$ws.Range("N105").Value = -5919

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 400564.72
$ws.Range("I4").Value = 385190.16
$ws.Range("K4").Value = 1155570.48
$ws.Range("M4").Value = -1155458.48
$ws.Range("H39").Value = 4706.3335
$ws.Range("J39").Value = 5004
$ws.Range("L39").Value = 15012
$ws.Range("N39").Value = -15600
$ws.Range("H55").Value = 1500
$ws.Range("I55").Value = 1428.5714
$ws.Range("J55").Value = 1666.6666
$ws.Range("K55").Value = 4285.7142
$ws.Range("L55").Value = 4999.9998
$ws.Range("M55").Value = -4108.7142
$ws.Range("N55").Value = -5353.9998
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 2333333
$ws.Range("I21").Value = 2333333
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 2333333
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -2333160
$ws.Range("N21").ClearContents()
$ws.Range("H30").Value = 2333333
$ws.Range("I30").Value = 2333333
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 2333333
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -2333228
$ws.Range("N30").ClearContents()
$ws.Range("H107").Value = 461
$ws.Range("I107").Value = 375.25
$ws.Range("J107").Value = 575.3333
$ws.Range("K107").Value = 375.25
$ws.Range("L107").Value = 575.3333
$ws.Range("M107").Value = 1544.75
$ws.Range("N107").Value = -4415.3333

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 21300
$ws.Range("I26").Value = 20000
$ws.Range("J26").Value = 22600
$ws.Range("K26").Value = 20000
$ws.Range("L26").Value = 22600
$ws.Range("M26").Value = -19705
$ws.Range("N26").Value = -23190

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 13621.417
$ws.Range("H84").Value = 13621.417
$ws.Range("H100").Value = 19224.5
$ws.Range("I100").Value = 19224.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 38449
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -37908
$ws.Range("N100").ClearContents()
